# Update the legacy GSC export data ("Chart" sheet):
#  - Drop the oldest date row ("2025-10-07") by deleting row 2; this shifts every
#    remaining row up by one and the unused shared string is dropped automatically.
#  - Append a new trailing row for the next day ("2026-01-05") with zeroed counts,
#    keeping the date as text (not an auto-converted date serial number).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

# Remove the oldest date row (row 2 holds "2025-10-07"); everything below moves up.
$ws.Rows.Item(2).Delete()

# Figure out where the new trailing row should go (right after the current last row).
$newRowIndex = $ws.UsedRange.Rows.Count + 1
$dateCell = $ws.Cells.Item($newRowIndex, 1)
$nonHttpsCell = $ws.Cells.Item($newRowIndex, 2)
$httpsCell = $ws.Cells.Item($newRowIndex, 3)

# Force the new date to be stored as text (like the existing date cells), not a
# date serial number: temporarily use a text number format while entering it...
$dateCell.NumberFormat = "@"
$dateCell.Formula = "2026-01-05"

# ...then restore the cell's formatting to match the rest of the date column by
# copying the (default) format from a sibling date cell, so no stray per-cell
# style lingers on the new cell.
$referenceDateCell = $ws.Cells.Item($newRowIndex - 1, 1)
$referenceDateCell.Copy()
$dateCell.PasteSpecial(-4122)
$excel.CutCopyMode = 0

$nonHttpsCell.Value = 0
$httpsCell.Value = 0
